$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 (I0) and J1 (IF) ---
# Copy formatting from the existing header cell H1 so the new headers
# match the bold/centered/bordered header style (style index 1).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data cells I2:J81 ---
$data = @(
    @(1,1),
    @(1,1),
    @(8,8),
    @(9,9),
    @(1,1),
    @(1,1),
    @(1,1),
    @(10,10),
    @(1,1),
    @(8,8),
    @(1,1),
    @(8,8),
    @(1,1),
    @(8,8),
    @(7,7),
    @(9,9),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(7,7),
    @(7,8),
    @(7,7),
    @(8,8),
    @(4,4),
    @(6,6),
    @(8,8),
    @(8,8),
    @(8,8),
    @(5,5),
    @(6,7),
    @(7,7),
    @(8,8),
    @(6,6),
    @(1,1),
    @(6,7),
    @(8,8),
    @(6,6),
    @(1,1),
    @(6,6),
    @(7,7),
    @(9,9),
    @(5,5),
    @(8,8),
    @(7,7),
    @(6,6),
    @(8,8),
    @(7,7),
    @(6,6),
    @(6,6),
    @(5,5),
    @(5,5),
    @(1,1),
    @(6,6),
    @(7,7),
    @(6,6),
    @(9,9),
    @(4,4),
    @(7,7),
    @(6,6),
    @(8,8),
    @(6,7),
    @(7,7),
    @(8,8),
    @(6,6),
    @(4,4),
    @(7,8),
    @(8,8),
    @(3,4),
    @(7,7),
    @(6,7),
    @(7,7),
    @(6,6),
    @(6,6),
    @(10,10),
    @(8,9),
    @(5,5),
    @(8,8),
    @(7,7)
)

$i = 0
foreach ($row in $data) {
    $r = $i + 2
    $ws.Cells.Item($r, 9).Value = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $i = $i + 1
}
